# Update column A values for specific rows to reflect updated KNN imputation
# results ("Update Name of Algo" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = -21.945
$ws.Range("A14").Value = -21.66
$ws.Range("A21").Value = -19.953
$ws.Range("A23").Value = -20.203
$ws.Range("A25").Value = -21.632
$ws.Range("A26").Value = -21.692
$ws.Range("A29").Value = -21.344
$ws.Range("A53").Value = -21.915
$ws.Range("A57").Value = -22.219
$ws.Range("A59").Value = -22.461
$ws.Range("A69").Value = -21.519
$ws.Range("A79").Value = -21.087
$ws.Range("A83").Value = -21.984
$ws.Range("A91").Value = -20.666
$ws.Range("A93").Value = -21.508
$ws.Range("A103").Value = -22.086
